# Riley Changes July 2nd
# The "Data" sheet had an extra/duplicate experimental run recorded in row 9
# (T (deg C)=90, AC=4.5%, but with stray D/T/Y values that don't belong to
# the 170/160/150/140/120-degree LSR sweep below it). Remove that row so the
# remaining rows of the sweep shift up into place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Make sure we're looking at the right sheet before editing, same as the
# workbook being saved with the Data tab active/selected.
$ws.Activate()

# Delete row 9 in its entirety - this shifts every row below it up by one,
# which is the only substantive change in the workbook.
$ws.Rows(9).EntireRow.Delete()

# Restore the cursor/selection to where it ended up after the edit.
$ws.Range("H23").Select()
